# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 42, pushing the existing
# rows 42-61 down to 43-62 (new dimension becomes A1:R62).
# The new row duplicates the market/product metadata of the old row 42
# but carries a new date (Fecha) and volume (Volumen).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 42..61 down to 43..62, leaving row 42 blank (except the
# inherited date-number format on column D).
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly record.
$ws.Range("A42").Value2 = 7
$ws.Range("B42").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C42").Value2 = "Ñuble"
$ws.Range("D42").Value2 = 44875
$ws.Range("E42").Value2 = 16
$ws.Range("F42").Value2 = 100112026
$ws.Range("G42").Value2 = "Haba"
$ws.Range("H42").Value2 = "Sin especificar"
$ws.Range("I42").Value2 = "Primera"
$ws.Range("J42").Value2 = 120
$ws.Range("K42").Value2 = 6500
$ws.Range("L42").Value2 = 7000
$ws.Range("M42").Value2 = 6750
$ws.Range("N42").Value2 = "$/saco 25 kilos"
$ws.Range("O42").Value2 = "Provincia de Diguillín"
$ws.Range("P42").Value2 = 270
$ws.Range("Q42").Value2 = 25
$ws.Range("R42").Value2 = "Hortaliza"
